$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 32500
$ws.Range("J93").Value = 32500
$ws.Range("L93").Value = 32500
$ws.Range("N93").Value = -37492
$ws.Range("H113").Value = 2275.875
$ws.Range("I113").Value = 2172.4285
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2172.4285
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1081.5715
$ws.Range("N113").Value = -9508
$ws.Range("H130").Value = 54656
$ws.Range("J130").Value = 54656
$ws.Range("L130").Value = 54656
$ws.Range("N130").Value = -64696

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1958.5
$ws.Range("I2").Value = 1837.3182
$ws.Range("J2").Value = 2625
$ws.Range("K2").Value = 1837.3182
$ws.Range("L2").Value = 2625
$ws.Range("M2").Value = -1724.3182
$ws.Range("N2").Value = -2851
$ws.Range("H45").Value = 2111.0557
$ws.Range("I45").Value = 1905.6666
$ws.Range("J45").Value = 3138
$ws.Range("K45").Value = 1905.6666
$ws.Range("L45").Value = 3138
$ws.Range("M45").Value = -1528.6666
$ws.Range("N45").Value = -3892
$ws.Range("H61").Value = 2098.1282
$ws.Range("I61").Value = 1942.4445
$ws.Range("J61").Value = 3966.3333
$ws.Range("K61").Value = 1942.4445
$ws.Range("L61").Value = 3966.3333
$ws.Range("M61").Value = -1730.4445
$ws.Range("N61").Value = -4390.3333
$ws.Range("H116").Value = 1958.5
$ws.Range("I116").Value = 1837.3182
$ws.Range("J116").Value = 2625
$ws.Range("K116").Value = 1837.3182
$ws.Range("L116").Value = 2625
$ws.Range("M116").Value = 456.6818000000001
$ws.Range("N116").Value = -7213
$ws.Range("H119").Value = 48979.2
$ws.Range("J119").Value = 48979.2
$ws.Range("L119").Value = 48979.2
$ws.Range("N119").Value = -58655.2
$ws.Range("H133").Value = 23545.234
$ws.Range("J133").Value = 23545.234
$ws.Range("L133").Value = 23545.234
$ws.Range("N133").Value = -28605.234
$ws.Range("H135").Value = 35784.23
$ws.Range("J135").Value = 35784.23
$ws.Range("L135").Value = 35784.23
$ws.Range("N135").Value = -45924.23
$ws.Range("H136").Value = 2098.1282
$ws.Range("I136").Value = 1942.4445
$ws.Range("J136").Value = 3966.3333
$ws.Range("K136").Value = 5827.333500000001
$ws.Range("L136").Value = 11898.9999
$ws.Range("M136").Value = -3277.333500000001
$ws.Range("N136").Value = -16998.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1958.5
$ws.Range("I3").Value = 1837.3182
$ws.Range("J3").Value = 2625
$ws.Range("K3").Value = 1837.3182
$ws.Range("L3").Value = 2625
$ws.Range("M3").Value = -1723.3182
$ws.Range("N3").Value = -2853
$ws.Range("H92").Value = 41700.5
$ws.Range("J92").Value = 41700.5
$ws.Range("L92").Value = 41700.5
$ws.Range("N92").Value = -46692.5
$ws.Range("H94").Value = 669.34485
$ws.Range("I94").Value = 589.13635
$ws.Range("J94").Value = 921.4286
$ws.Range("K94").Value = 589.13635
$ws.Range("L94").Value = 921.4286
$ws.Range("M94").Value = -138.13635
$ws.Range("N94").Value = -1823.4286
$ws.Range("H126").Value = 48181.332
$ws.Range("J126").Value = 48181.332
$ws.Range("L126").Value = 48181.332
$ws.Range("N126").Value = -58061.332
$ws.Range("H134").Value = 2961.8
$ws.Range("I134").Value = 2580.5
$ws.Range("K134").Value = 7741.5
$ws.Range("M134").Value = -5206.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 27973.111
$ws.Range("J41").Value = 33965.43
$ws.Range("L41").Value = 33965.43
$ws.Range("N41").Value = -34821.43
$ws.Range("H50").Value = 37963
$ws.Range("J50").Value = 37963
$ws.Range("L50").Value = 37963
$ws.Range("N50").Value = -39213
$ws.Range("H59").Value = 34526.668
$ws.Range("J59").Value = 34832
$ws.Range("L59").Value = 34832
$ws.Range("N59").Value = -37122
$ws.Range("H60").Value = 23691.234
$ws.Range("J60").Value = 23691.234
$ws.Range("L60").Value = 23691.234
$ws.Range("N60").Value = -24713.234
$ws.Range("H68").Value = 165000
$ws.Range("J68").Value = 165000
$ws.Range("L68").Value = 165000
$ws.Range("N68").Value = -166498
$ws.Range("H70").Value = 31168.857
$ws.Range("J70").Value = 31168.857
$ws.Range("L70").Value = 31168.857
$ws.Range("N70").Value = -31798.857
$ws.Range("H71").Value = 165000
$ws.Range("J71").Value = 165000
$ws.Range("L71").Value = 495000
$ws.Range("N71").Value = -502488
$ws.Range("H73").Value = 31168.857
$ws.Range("J73").Value = 31168.857
$ws.Range("L73").Value = 31168.857
$ws.Range("N73").Value = -33352.857
$ws.Range("H80").Value = 28844.125
$ws.Range("J80").Value = 28844.125
$ws.Range("L80").Value = 28844.125
$ws.Range("N80").Value = -31090.125
$ws.Range("H83").Value = 28844.125
$ws.Range("J83").Value = 28844.125
$ws.Range("L83").Value = 86532.375
$ws.Range("N83").Value = -97764.375
$ws.Range("H92").Value = 38400.668
$ws.Range("J92").Value = 38400.668
$ws.Range("L92").Value = 38400.668
$ws.Range("N92").Value = -43392.668
$ws.Range("H99").Value = 2406
$ws.Range("I99").Value = 2338.8572
$ws.Range("K99").Value = 2338.8572
$ws.Range("M99").Value = -840.8571999999999
$ws.Range("H126").Value = 2406
$ws.Range("I126").Value = 2338.8572
$ws.Range("K126").Value = 7016.571599999999
$ws.Range("M126").Value = -4546.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19400
$ws.Range("J15").Value = 19400
$ws.Range("L15").Value = 19400
$ws.Range("N15").Value = -19976
$ws.Range("H81").Value = 19400
$ws.Range("J81").Value = 19400
$ws.Range("L81").Value = 19400
$ws.Range("N81").Value = -21396
$ws.Range("H84").Value = 19400
$ws.Range("J84").Value = 19400
$ws.Range("L84").Value = 58200
$ws.Range("N84").Value = -68184
$ws.Range("H123").Value = 24163
$ws.Range("J123").Value = 24163
$ws.Range("L123").Value = 24163
$ws.Range("N123").Value = -29063
$ws.Range("H130").Value = 49598.4
$ws.Range("J130").Value = 49598.4
$ws.Range("L130").Value = 49598.4
$ws.Range("N130").Value = -59638.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 45433
$ws.Range("J119").Value = 45433
$ws.Range("L119").Value = 45433
$ws.Range("N119").Value = -55109

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 29374
$ws.Range("J93").Value = 29374
$ws.Range("L93").Value = 29374
$ws.Range("N93").Value = -34366
$ws.Range("H98").Value = 45301
$ws.Range("J98").Value = 45301
$ws.Range("L98").Value = 45301
$ws.Range("N98").Value = -51291
$ws.Range("H119").Value = 47996
$ws.Range("J119").Value = 47996
$ws.Range("L119").Value = 47996
$ws.Range("N119").Value = -57672
$ws.Range("H120").Value = 43110.4
$ws.Range("J120").Value = 43110.4
$ws.Range("L120").Value = 43110.4
$ws.Range("N120").Value = -52786.4
$ws.Range("H122").Value = 58731628
$ws.Range("I122").Value = 105715210
$ws.Range("J122").Value = 2157.5
$ws.Range("K122").Value = 317145630
$ws.Range("L122").Value = 6472.5
$ws.Range("M122").Value = -317143180
$ws.Range("N122").Value = -11372.5
